$d = $word.ActiveDocument

# Locate the two paragraphs that need restructuring:
#   P1: "Ugyldig dato består af dato der ligger i datid på det tidspunkt."
#   P2: "Ekstra tilbehørs pris, der afhængig af ikke-obligatoriske oplysninger bruges for beregningen. "
$count = $d.Paragraphs.Count
$idxUgyldig = -1
$idxEkstra = -1
for ($i = 1; $i -le $count; $i++) {
    $t = $d.Paragraphs.Item($i).Range.Text
    if ($t -like "*Ugyldig dato*tidspunkt*") {
        $idxUgyldig = $i
    }
    if ($t -like "*Ekstra tilbehør*beregningen*") {
        $idxEkstra = $i
    }
}

$pUgyldig = $d.Paragraphs.Item($idxUgyldig)
$pEkstra = $d.Paragraphs.Item($idxEkstra)

# Range spanning both paragraphs in full (start of first through end of second,
# including both paragraph marks).
$r = $d.Range($pUgyldig.Range.Start, $pEkstra.Range.End)

# Replace that span with the target markup: the trailing "_GoBack" bookmark now
# sits at the end of the first paragraph (right after "tidspunkt.") instead of
# splitting " bruges" / " for beregningen" in the second paragraph, so those two
# runs merge back into a single " bruges for beregningen" run.
$xmlFrag = '<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p w:rsidR="001824A8" w:rsidRPr="001824A8" w:rsidRDefault="001824A8" w:rsidP="001824A8"><w:pPr><w:rPr><w:rFonts w:asciiTheme="minorHAnsi" w:hAnsiTheme="minorHAnsi"/><w:color w:val="FF0000"/><w:sz w:val="22"/><w:szCs w:val="22"/><w:lang w:val="da-DK"/></w:rPr></w:pPr><w:r w:rsidRPr="001824A8"><w:rPr><w:rFonts w:asciiTheme="minorHAnsi" w:hAnsiTheme="minorHAnsi"/><w:color w:val="FF0000"/><w:sz w:val="22"/><w:szCs w:val="22"/><w:lang w:val="da-DK"/></w:rPr><w:t>Ugyldig dato består af dato der ligger i datid på det tidspunkt.</w:t></w:r><w:bookmarkStart w:id="0" w:name="_GoBack"/><w:bookmarkEnd w:id="0"/></w:p><w:p w:rsidR="00835A08" w:rsidRDefault="00835A08" w:rsidP="00835A08"><w:pPr><w:pStyle w:val="BodyA"/><w:rPr><w:rFonts w:asciiTheme="minorHAnsi" w:hAnsiTheme="minorHAnsi"/></w:rPr></w:pPr><w:r w:rsidRPr="00835A08"><w:rPr><w:rFonts w:asciiTheme="minorHAnsi" w:hAnsiTheme="minorHAnsi"/><w:color w:val="FF0000"/></w:rPr><w:t>Ekstra tilbehørs pris</w:t></w:r><w:r w:rsidR="00530681"><w:rPr><w:rFonts w:asciiTheme="minorHAnsi" w:hAnsiTheme="minorHAnsi"/><w:color w:val="FF0000"/></w:rPr><w:t>, der afhængig af ikke-obligatoriske oplysninger</w:t></w:r><w:r w:rsidRPr="00835A08"><w:rPr><w:rFonts w:asciiTheme="minorHAnsi" w:hAnsiTheme="minorHAnsi"/><w:color w:val="FF0000"/></w:rPr><w:t xml:space="preserve"> bruges for beregningen</w:t></w:r><w:r><w:rPr><w:rFonts w:asciiTheme="minorHAnsi" w:hAnsiTheme="minorHAnsi"/></w:rPr><w:t xml:space="preserve">. </w:t></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'

$r.InsertXML($xmlFrag)
